$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "kk"
$ws.Range("B8").Value = 33

$ws.Range("B9").Select()
